# Fill in the match results (goals, result/winner, and points) for the
# set of Matchweek 9/10 fixtures that were previously blank.
#
# Columns: A=MatchWeek, B=Home, C=Away, D=HomeGoals, E=AwayGoals,
#          F=Result (winner team name, or "Draw"), N=HomePoints, O=AwayPoints

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$matches = @(
    @{ Row = 14;  MatchWeek = 9;  HomeGoals = 1; AwayGoals = 0; Result = "Arsenal";            HomePoints = 3; AwayPoints = 0 },
    @{ Row = 33;  MatchWeek = 10; HomeGoals = 0; AwayGoals = 2; Result = "Arsenal";            HomePoints = 0; AwayPoints = 3 },
    @{ Row = 42;  MatchWeek = 10; HomeGoals = 3; AwayGoals = 1; Result = "Manchester City";    HomePoints = 3; AwayPoints = 0 },
    @{ Row = 57;  MatchWeek = 9;  HomeGoals = 1; AwayGoals = 0; Result = "Aston Villa";         HomePoints = 3; AwayPoints = 0 },
    @{ Row = 90;  MatchWeek = 10; HomeGoals = 2; AwayGoals = 0; Result = "Liverpool";           HomePoints = 3; AwayPoints = 0 },
    @{ Row = 101; MatchWeek = 9;  HomeGoals = 3; AwayGoals = 2; Result = "Brentford";           HomePoints = 3; AwayPoints = 0 },
    @{ Row = 136; MatchWeek = 9;  HomeGoals = 2; AwayGoals = 0; Result = "AFC Bournemouth";     HomePoints = 3; AwayPoints = 0 },
    @{ Row = 143; MatchWeek = 10; HomeGoals = 0; AwayGoals = 1; Result = "Chelsea";             HomePoints = 0; AwayPoints = 3 },
    @{ Row = 144; MatchWeek = 9;  HomeGoals = 1; AwayGoals = 2; Result = "Sunderland";          HomePoints = 0; AwayPoints = 3 },
    @{ Row = 183; MatchWeek = 9;  HomeGoals = 0; AwayGoals = 3; Result = "Tottenham Hotspur";   HomePoints = 0; AwayPoints = 3 },
    @{ Row = 208; MatchWeek = 10; HomeGoals = 1; AwayGoals = 1; Result = "Draw";                HomePoints = 1; AwayPoints = 1 },
    @{ Row = 240; MatchWeek = 10; HomeGoals = 2; AwayGoals = 0; Result = "Crystal Palace";      HomePoints = 3; AwayPoints = 0 },
    @{ Row = 250; MatchWeek = 9;  HomeGoals = 4; AwayGoals = 2; Result = "Manchester United";   HomePoints = 3; AwayPoints = 0 },
    @{ Row = 267; MatchWeek = 10; HomeGoals = 2; AwayGoals = 2; Result = "Draw";                HomePoints = 1; AwayPoints = 1 },
    @{ Row = 280; MatchWeek = 10; HomeGoals = 3; AwayGoals = 0; Result = "Brighton & Hove Albion"; HomePoints = 3; AwayPoints = 0 },
    @{ Row = 326; MatchWeek = 9;  HomeGoals = 2; AwayGoals = 1; Result = "Newcastle United";    HomePoints = 3; AwayPoints = 0 },
    @{ Row = 337; MatchWeek = 10; HomeGoals = 3; AwayGoals = 1; Result = "West Ham United";     HomePoints = 3; AwayPoints = 0 },
    @{ Row = 350; MatchWeek = 10; HomeGoals = 3; AwayGoals = 0; Result = "Fulham";              HomePoints = 3; AwayPoints = 0 },
    @{ Row = 358; MatchWeek = 9;  HomeGoals = 2; AwayGoals = 1; Result = "Leeds United";        HomePoints = 3; AwayPoints = 0 },
    @{ Row = 375; MatchWeek = 9;  HomeGoals = 2; AwayGoals = 3; Result = "Burnley";             HomePoints = 0; AwayPoints = 3 }
)

foreach ($m in $matches) {
    $r = $m.Row
    $ws.Cells.Item($r, 1).Value = $m.MatchWeek   # A - MatchWeek
    $ws.Cells.Item($r, 4).Value = $m.HomeGoals   # D - HomeGoals
    $ws.Cells.Item($r, 5).Value = $m.AwayGoals   # E - AwayGoals
    $ws.Cells.Item($r, 6).Value = $m.Result      # F - Result
    $ws.Cells.Item($r, 14).Value = $m.HomePoints # N - HomePoints
    $ws.Cells.Item($r, 15).Value = $m.AwayPoints # O - AwayPoints
}
